# Trade #74 closed at 2026-02-17 15:48:55 - unknown UNKNOWN +0.000%
#
# Updates the "Summary" and "Strategy Status" roll-up figures to reflect the
# newly closed trade, and appends the new trade record (row 75, Trade # 74)
# to both the "All Trades" and "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.15   # Current Capital
$summary.Range("B4").Value = 0.14      # Total P&L $
$summary.Range("B5").Value = 0.04      # Total P&L %
$summary.Range("B6").Value = 74        # Total Trades
$summary.Range("B7").Value = 24        # Winning Trades
$summary.Range("B9").Value = 32.43     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.15     # Capital
$status.Range("D4").Value = 74         # Trades
$status.Range("E4").Value = 0.14       # P&L $
$status.Range("F4").Value = 0.15       # P&L %
$status.Range("G4").Value = 32.43      # Win Rate %

# ---------------------------------------------------------------------
# New trade row (row 75) appended to "All Trades" and "MarketMaking"
# ---------------------------------------------------------------------
# Columns: A Trade#, B Date, C Time, D Strategy, E Side, F Entry Price,
#          G Exit Price, H Status, I P&L %, J P&L $, K Capital After,
#          L Entry Slippage, M Exit Slippage, N Confidence,
#          O Entry Reason, P Exit Reason, Q Duration (min)
$newRowNum = 75
$newRowValues = @{
    1  = 74
    2  = "2026-02-17"
    3  = "15:48:48"
    4  = "MarketMaking"
    5  = "UP"
    6  = 0.68
    7  = 0.71
    8  = "CLOSED"
    9  = 4.4118
    10 = 0.03
    11 = 100.15
    12 = 0
    13 = 0
    14 = 0.6
    15 = "Normal spread capture: 19600 bps"
    16 = "early_exit"
    17 = 0.13
}
# Columns B (date) and C (time) look like dates/times to Excel's type
# inference, so they must be pinned to Text before the assignment (and the
# style reset back to Normal afterwards) to keep them as literal strings,
# matching the rest of the column which stores these as plain text.
$textCols = @(2, 3)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($col in $newRowValues.Keys) {
        $cell = $ws.Cells.Item($newRowNum, $col)
        if ($textCols -contains $col) {
            $cell.NumberFormat = "@"
            $cell.Value = $newRowValues[$col]
            $cell.Style = "Normal"
        } else {
            $cell.Value = $newRowValues[$col]
        }
    }
}
